$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 946.7692
$ws.Range("I92").Value = 993.5833
$ws.Range("J92").Value = 385
$ws.Range("K92").Value = 993.5833
$ws.Range("L92").Value = 385
$ws.Range("M92").Value = 254.4167
$ws.Range("N92").Value = -2881
$ws.Range("H101").Value = 11907594
$ws.Range("I101").Value = 23811484
$ws.Range("J101").Value = 3704.8333
$ws.Range("K101").Value = 71434452
$ws.Range("L101").Value = 11114.4999
$ws.Range("M101").Value = -71432830
$ws.Range("N101").Value = -14358.4999
$ws.Range("H135").Value = 5527.6294
$ws.Range("I135").Value = 7037.0557
$ws.Range("J135").Value = 2508.7778
$ws.Range("K135").Value = 63333.5013
$ws.Range("L135").Value = 22579.0002
$ws.Range("M135").Value = -60798.5013
$ws.Range("N135").Value = -27649.0002
$ws.Range("H138").Value = 2283.4285
$ws.Range("J138").Value = 4084.5
$ws.Range("L138").Value = 12253.5
$ws.Range("N138").Value = -22533.5
$ws.Range("H141").Value = 7084.1055
$ws.Range("J141").Value = 5749.75
$ws.Range("L141").Value = 17249.25
$ws.Range("N141").Value = -27609.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6945.846
$ws.Range("I32").Value = 6885.961
$ws.Range("K32").Value = 6885.961
$ws.Range("M32").Value = -6598.961
$ws.Range("H45").Value = 59810.918
$ws.Range("I45").Value = 115287.78
$ws.Range("J45").Value = 4334.0557
$ws.Range("K45").Value = 115287.78
$ws.Range("L45").Value = 4334.0557
$ws.Range("M45").Value = -114910.78
$ws.Range("N45").Value = -5088.0557
$ws.Range("H122").Value = 1505995.9
$ws.Range("I122").Value = 6423.1
$ws.Range("J122").Value = 3005568.8
$ws.Range("K122").Value = 19269.3
$ws.Range("L122").Value = 9016706.399999999
$ws.Range("M122").Value = -16819.3
$ws.Range("N122").Value = -9021606.399999999
$ws.Range("H132").Value = 3152.1516
$ws.Range("I132").Value = 3115.074
$ws.Range("K132").Value = 9345.222
$ws.Range("M132").Value = -6815.222

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7576.2095
$ws.Range("I94").Value = 9494.031999999999
$ws.Range("J94").Value = 2621.8333
$ws.Range("K94").Value = 9494.031999999999
$ws.Range("L94").Value = 2621.8333
$ws.Range("M94").Value = -9043.031999999999
$ws.Range("N94").Value = -3523.8333
$ws.Range("H134").Value = 7289.0435
$ws.Range("I134").Value = 8351.210999999999
$ws.Range("K134").Value = 25053.633
$ws.Range("M134").Value = -22518.633

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9298.85
$ws.Range("I31").Value = 14000.8
$ws.Range("J31").Value = 4596.9
$ws.Range("K31").Value = 14000.8
$ws.Range("L31").Value = 4596.9
$ws.Range("M31").Value = -13705.8
$ws.Range("N31").Value = -5186.9
$ws.Range("H34").Value = 9298.85
$ws.Range("I34").Value = 14000.8
$ws.Range("J34").Value = 4596.9
$ws.Range("K34").Value = 14000.8
$ws.Range("L34").Value = 4596.9
$ws.Range("M34").Value = -13798.8
$ws.Range("N34").Value = -5000.9
$ws.Range("H58").Value = 2456
$ws.Range("I58").Value = 2342.0715
$ws.Range("J58").Value = 2746
$ws.Range("K58").Value = 2342.0715
$ws.Range("L58").Value = 2746
$ws.Range("M58").Value = -2139.0715
$ws.Range("N58").Value = -3152
$ws.Range("H131").Value = 44444
$ws.Range("J131").Value = 44444
$ws.Range("L131").Value = 44444
$ws.Range("N131").Value = -54524
$ws.Range("H132").Value = 1677.375
$ws.Range("I132").Value = 1631.2858
$ws.Range("K132").Value = 4893.857400000001
$ws.Range("M132").Value = -2363.857400000001
$ws.Range("H134").Value = 6416.6665
$ws.Range("I134").Value = 8319.611000000001
$ws.Range("J134").Value = 2610.7778
$ws.Range("K134").Value = 24958.833
$ws.Range("L134").Value = 7832.3334
$ws.Range("M134").Value = -22423.833
$ws.Range("N134").Value = -12902.3334
$ws.Range("H136").Value = 2456
$ws.Range("I136").Value = 2342.0715
$ws.Range("J136").Value = 2746
$ws.Range("K136").Value = 7026.2145
$ws.Range("L136").Value = 8238
$ws.Range("M136").Value = -4476.2145
$ws.Range("N136").Value = -13338
$ws.Range("H141").Value = 330717.47
$ws.Range("J141").Value = 377211.53
$ws.Range("L141").Value = 377211.53
$ws.Range("N141").Value = -387571.53

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 383.625
$ws.Range("I92").Value = 361.5
$ws.Range("J92").Value = 450
$ws.Range("K92").Value = 1084.5
$ws.Range("L92").Value = 1350
$ws.Range("M92").Value = 163.5
$ws.Range("N92").Value = -3846
$ws.Range("H98").Value = 1927.9
$ws.Range("J98").Value = 1595.6666
$ws.Range("L98").Value = 4786.9998
$ws.Range("N98").Value = -7782.9998
$ws.Range("H131").Value = 2752.0762
$ws.Range("J131").Value = 1972.2235
$ws.Range("L131").Value = 5916.6705
$ws.Range("N131").Value = -15996.6705

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 48998.5
$ws.Range("J46").Value = 48998.5
$ws.Range("L46").Value = 48998.5
$ws.Range("N46").Value = -49310.5
$ws.Range("H126").Value = 6627
$ws.Range("I126").Value = 12985.637
$ws.Range("J126").Value = 2945.6843
$ws.Range("K126").Value = 38956.911
$ws.Range("L126").Value = 8837.052899999999
$ws.Range("M126").Value = -36486.911
$ws.Range("N126").Value = -13777.0529
$ws.Range("H132").Value = 3927.8958
$ws.Range("I132").Value = 4345.6113
$ws.Range("J132").Value = 2674.75
$ws.Range("K132").Value = 13036.8339
$ws.Range("L132").Value = 8024.25
$ws.Range("M132").Value = -10506.8339
$ws.Range("N132").Value = -13084.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17172.516
$ws.Range("I7").Value = 24266.8
$ws.Range("K7").Value = 24266.8
$ws.Range("M7").Value = -24154.8
$ws.Range("H22").Value = 12016.667
$ws.Range("I22").Value = 22677.889
$ws.Range("K22").Value = 22677.889
$ws.Range("M22").Value = -22382.889
$ws.Range("H27").Value = 12016.667
$ws.Range("I27").Value = 22677.889
$ws.Range("K27").Value = 22677.889
$ws.Range("M27").Value = -22570.889
$ws.Range("H122").Value = 4445.9756
$ws.Range("I122").Value = 4146.9688
$ws.Range("K122").Value = 12440.9064
$ws.Range("M122").Value = -9990.9064
$ws.Range("H126").Value = 17172.516
$ws.Range("I126").Value = 24266.8
$ws.Range("K126").Value = 72800.39999999999
$ws.Range("M126").Value = -70330.39999999999
$ws.Range("H132").Value = 786442.5
$ws.Range("I132").Value = 994901.2
$ws.Range("K132").Value = 2984703.6
$ws.Range("M132").Value = -2982173.6
$ws.Range("H136").Value = 3705.3948
$ws.Range("I136").Value = 2821.0344
$ws.Range("K136").Value = 8463.1032
$ws.Range("M136").Value = -5913.1032

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 18984.38
$ws.Range("I100").Value = 4834.25
$ws.Range("J100").Value = 64264.8
$ws.Range("K100").Value = 9668.5
$ws.Range("L100").Value = 128529.6
$ws.Range("M100").Value = -9127.5
$ws.Range("N100").Value = -129611.6
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180
$ws.Range("H126").Value = 33241.54
$ws.Range("I126").Value = 41794.8
$ws.Range("K126").Value = 125384.4
$ws.Range("M126").Value = -122914.4
$ws.Range("H132").Value = 6144.75
$ws.Range("I132").Value = 7148.3193
$ws.Range("J132").Value = 3370.1765
$ws.Range("K132").Value = 21444.9579
$ws.Range("L132").Value = 10110.5295
$ws.Range("M132").Value = -18914.9579
$ws.Range("N132").Value = -15170.5295
